# Reroute test data update
# - Append ", FragilePAK" to the AvailableCarriers list shown in C2
# - Move the active selection from B2 to C2

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = "ABF Freight, BTX Global Logistics, Ceva,FC Test Carrier, UPS, FragilePAK"

$ws.Range("C2").Select()
